$d = $word.ActiveDocument

$d.Content.Find.Execute("311÷2=155, 1", $true, $false, $false, $false, $false, $true, 1, $false, "397÷2=198, 1", 2) | Out-Null
$d.Content.Find.Execute("355÷4=88, 3", $true, $false, $false, $false, $false, $true, 1, $false, "564÷2=282, 0", 2) | Out-Null
$d.Content.Find.Execute("393÷6=65, 3", $true, $false, $false, $false, $false, $true, 1, $false, "514÷2=257, 0", 2) | Out-Null
$d.Content.Find.Execute("408÷4=102, 0", $true, $false, $false, $false, $false, $true, 1, $false, "962÷9=106, 8", 2) | Out-Null
$d.Content.Find.Execute("425÷2=212, 1", $true, $false, $false, $false, $false, $true, 1, $false, "518÷9=57, 5", 2) | Out-Null
$d.Content.Find.Execute("993÷5=198, 3", $true, $false, $false, $false, $false, $true, 1, $false, "703÷5=140, 3", 2) | Out-Null
$d.Content.Find.Execute("489÷8=61, 1", $true, $false, $false, $false, $false, $true, 1, $false, "294÷5=58, 4", 2) | Out-Null
$d.Content.Find.Execute("499÷9=55, 4", $true, $false, $false, $false, $false, $true, 1, $false, "715÷8=89, 3", 2) | Out-Null
$d.Content.Find.Execute("683÷3=227, 2", $true, $false, $false, $false, $false, $true, 1, $false, "366÷3=122, 0", 2) | Out-Null
$d.Content.Find.Execute("475÷4=118, 3", $true, $false, $false, $false, $false, $true, 1, $false, "333÷9=37, 0", 2) | Out-Null
$d.Content.Find.Execute("885÷2=442, 1", $true, $false, $false, $false, $false, $true, 1, $false, "470÷4=117, 2", 2) | Out-Null
$d.Content.Find.Execute("446÷4=111, 2", $true, $false, $false, $false, $false, $true, 1, $false, "723÷3=241, 0", 2) | Out-Null
$d.Content.Find.Execute("822÷2=411, 0", $true, $false, $false, $false, $false, $true, 1, $false, "128÷8=16, 0", 2) | Out-Null
$d.Content.Find.Execute("727÷9=80, 7", $true, $false, $false, $false, $false, $true, 1, $false, "717÷6=119, 3", 2) | Out-Null
$d.Content.Find.Execute("526÷8=65, 6", $true, $false, $false, $false, $false, $true, 1, $false, "639÷5=127, 4", 2) | Out-Null
$d.Content.Find.Execute("290÷2=145, 0", $true, $false, $false, $false, $false, $true, 1, $false, "447÷5=89, 2", 2) | Out-Null
$d.Content.Find.Execute("533÷7=76, 1", $true, $false, $false, $false, $false, $true, 1, $false, "898÷3=299, 1", 2) | Out-Null
$d.Content.Find.Execute("411÷8=51, 3", $true, $false, $false, $false, $false, $true, 1, $false, "159÷8=19, 7", 2) | Out-Null
$d.Content.Find.Execute("546÷9=60, 6", $true, $false, $false, $false, $false, $true, 1, $false, "960÷9=106, 6", 2) | Out-Null
$d.Content.Find.Execute("975÷8=121, 7", $true, $false, $false, $false, $false, $true, 1, $false, "471÷8=58, 7", 2) | Out-Null
$d.Content.Find.Execute("621÷6=103, 3", $true, $false, $false, $false, $false, $true, 1, $false, "444÷2=222, 0", 2) | Out-Null
$d.Content.Find.Execute("316÷2=158, 0", $true, $false, $false, $false, $false, $true, 1, $false, "348÷9=38, 6", 2) | Out-Null
$d.Content.Find.Execute("579÷7=82, 5", $true, $false, $false, $false, $false, $true, 1, $false, "497÷9=55, 2", 2) | Out-Null
$d.Content.Find.Execute("274÷4=68, 2", $true, $false, $false, $false, $false, $true, 1, $false, "326÷5=65, 1", 2) | Out-Null
$d.Content.Find.Execute("561÷8=70, 1", $true, $false, $false, $false, $false, $true, 1, $false, "729÷9=81, 0", 2) | Out-Null
